# Applies the weekly data refresh for the Pepino dulce (Agro Chillan)
# subset: dates, volumes, prices and quality grades for rows 2-37 are
# shuffled to reflect the new weekly sample.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = 44839
$ws.Cells.Item(2, 10).Value2 = 100
$ws.Cells.Item(2, 11).Value2 = 17000
$ws.Cells.Item(2, 12).Value2 = 18000
$ws.Cells.Item(2, 13).Value2 = 17500
$ws.Cells.Item(2, 16).Value2 = 972
$ws.Cells.Item(3, 4).Value2 = 44832
$ws.Cells.Item(3, 10).Value2 = 60
$ws.Cells.Item(4, 4).Value2 = 44651
$ws.Cells.Item(4, 10).Value2 = 60
$ws.Cells.Item(4, 11).Value2 = 15000
$ws.Cells.Item(4, 12).Value2 = 16000
$ws.Cells.Item(4, 13).Value2 = 15500
$ws.Cells.Item(4, 16).Value2 = 861
$ws.Cells.Item(5, 4).Value2 = 44635
$ws.Cells.Item(5, 11).Value2 = 15000
$ws.Cells.Item(5, 12).Value2 = 16000
$ws.Cells.Item(5, 13).Value2 = 15500
$ws.Cells.Item(5, 16).Value2 = 861
$ws.Cells.Item(6, 4).Value2 = 44637
$ws.Cells.Item(6, 10).Value2 = 100
$ws.Cells.Item(7, 4).Value2 = 44799
$ws.Cells.Item(7, 10).Value2 = 60
$ws.Cells.Item(7, 11).Value2 = 15000
$ws.Cells.Item(7, 12).Value2 = 16000
$ws.Cells.Item(7, 13).Value2 = 15500
$ws.Cells.Item(7, 16).Value2 = 861
$ws.Cells.Item(8, 4).Value2 = 44642
$ws.Cells.Item(8, 10).Value2 = 100
$ws.Cells.Item(8, 11).Value2 = 15000
$ws.Cells.Item(8, 12).Value2 = 16000
$ws.Cells.Item(8, 13).Value2 = 15500
$ws.Cells.Item(8, 16).Value2 = 861
$ws.Cells.Item(9, 4).Value2 = 44790
$ws.Cells.Item(9, 10).Value2 = 60
$ws.Cells.Item(9, 11).Value2 = 17000
$ws.Cells.Item(9, 12).Value2 = 18000
$ws.Cells.Item(9, 13).Value2 = 17500
$ws.Cells.Item(9, 16).Value2 = 972
$ws.Cells.Item(11, 4).Value2 = 44830
$ws.Cells.Item(11, 11).Value2 = 17000
$ws.Cells.Item(11, 12).Value2 = 17000
$ws.Cells.Item(11, 13).Value2 = 17000
$ws.Cells.Item(11, 16).Value2 = 944
$ws.Cells.Item(12, 4).Value2 = 44664
$ws.Cells.Item(12, 10).Value2 = 160
$ws.Cells.Item(12, 11).Value2 = 15000
$ws.Cells.Item(12, 12).Value2 = 16000
$ws.Cells.Item(12, 13).Value2 = 15500
$ws.Cells.Item(12, 16).Value2 = 861
$ws.Cells.Item(13, 4).Value2 = 44769
$ws.Cells.Item(13, 11).Value2 = 17000
$ws.Cells.Item(13, 12).Value2 = 18000
$ws.Cells.Item(13, 13).Value2 = 17500
$ws.Cells.Item(13, 16).Value2 = 972
$ws.Cells.Item(14, 4).Value2 = 44785
$ws.Cells.Item(14, 10).Value2 = 80
$ws.Cells.Item(15, 4).Value2 = 44818
$ws.Cells.Item(15, 11).Value2 = 15000
$ws.Cells.Item(15, 12).Value2 = 15000
$ws.Cells.Item(15, 13).Value2 = 15000
$ws.Cells.Item(15, 16).Value2 = 833
$ws.Cells.Item(16, 4).Value2 = 44804
$ws.Cells.Item(16, 10).Value2 = 100
$ws.Cells.Item(16, 11).Value2 = 15000
$ws.Cells.Item(16, 12).Value2 = 16000
$ws.Cells.Item(16, 13).Value2 = 15500
$ws.Cells.Item(16, 16).Value2 = 861
$ws.Cells.Item(17, 4).Value2 = 44659
$ws.Cells.Item(18, 4).Value2 = 44819
$ws.Cells.Item(18, 10).Value2 = 60
$ws.Cells.Item(18, 11).Value2 = 15000
$ws.Cells.Item(18, 12).Value2 = 15000
$ws.Cells.Item(18, 13).Value2 = 15000
$ws.Cells.Item(18, 16).Value2 = 833
$ws.Cells.Item(19, 4).Value2 = 44384
$ws.Cells.Item(19, 10).Value2 = 120
$ws.Cells.Item(19, 11).Value2 = 17000
$ws.Cells.Item(19, 12).Value2 = 18000
$ws.Cells.Item(19, 13).Value2 = 17500
$ws.Cells.Item(19, 16).Value2 = 972
$ws.Cells.Item(20, 4).Value2 = 44384
$ws.Cells.Item(20, 9).Value2 = "Segunda"
$ws.Cells.Item(20, 10).Value2 = 60
$ws.Cells.Item(20, 12).Value2 = 15000
$ws.Cells.Item(20, 13).Value2 = 15000
$ws.Cells.Item(20, 16).Value2 = 833
$ws.Cells.Item(21, 4).Value2 = 44658
$ws.Cells.Item(21, 10).Value2 = 80
$ws.Cells.Item(22, 4).Value2 = 44645
$ws.Cells.Item(22, 10).Value2 = 60
$ws.Cells.Item(22, 11).Value2 = 15000
$ws.Cells.Item(22, 12).Value2 = 16000
$ws.Cells.Item(22, 13).Value2 = 15500
$ws.Cells.Item(22, 16).Value2 = 861
$ws.Cells.Item(23, 4).Value2 = 44813
$ws.Cells.Item(23, 10).Value2 = 100
$ws.Cells.Item(23, 11).Value2 = 14000
$ws.Cells.Item(23, 12).Value2 = 15000
$ws.Cells.Item(23, 13).Value2 = 14500
$ws.Cells.Item(23, 16).Value2 = 806
$ws.Cells.Item(24, 4).Value2 = 44809
$ws.Cells.Item(24, 11).Value2 = 14000
$ws.Cells.Item(24, 12).Value2 = 15000
$ws.Cells.Item(24, 13).Value2 = 14500
$ws.Cells.Item(24, 16).Value2 = 806
$ws.Cells.Item(25, 4).Value2 = 44791
$ws.Cells.Item(25, 10).Value2 = 80
$ws.Cells.Item(25, 11).Value2 = 17000
$ws.Cells.Item(25, 12).Value2 = 18000
$ws.Cells.Item(25, 13).Value2 = 17500
$ws.Cells.Item(25, 16).Value2 = 972
$ws.Cells.Item(26, 4).Value2 = 44797
$ws.Cells.Item(26, 10).Value2 = 80
$ws.Cells.Item(26, 11).Value2 = 16000
$ws.Cells.Item(26, 12).Value2 = 17000
$ws.Cells.Item(26, 13).Value2 = 16500
$ws.Cells.Item(26, 16).Value2 = 917
$ws.Cells.Item(27, 4).Value2 = 44775
$ws.Cells.Item(27, 10).Value2 = 100
$ws.Cells.Item(28, 4).Value2 = 44847
$ws.Cells.Item(28, 10).Value2 = 120
$ws.Cells.Item(28, 12).Value2 = 17000
$ws.Cells.Item(28, 13).Value2 = 17000
$ws.Cells.Item(28, 16).Value2 = 944
$ws.Cells.Item(29, 4).Value2 = 44628
$ws.Cells.Item(29, 10).Value2 = 60
$ws.Cells.Item(30, 4).Value2 = 44811
$ws.Cells.Item(30, 10).Value2 = 60
$ws.Cells.Item(30, 11).Value2 = 14000
$ws.Cells.Item(30, 12).Value2 = 15000
$ws.Cells.Item(30, 13).Value2 = 14500
$ws.Cells.Item(30, 16).Value2 = 806
$ws.Cells.Item(31, 4).Value2 = 44782
$ws.Cells.Item(31, 9).Value2 = "Primera"
$ws.Cells.Item(31, 10).Value2 = 120
$ws.Cells.Item(31, 11).Value2 = 17000
$ws.Cells.Item(31, 12).Value2 = 18000
$ws.Cells.Item(31, 13).Value2 = 17500
$ws.Cells.Item(31, 16).Value2 = 972
$ws.Cells.Item(32, 4).Value2 = 44649
$ws.Cells.Item(33, 4).Value2 = 44771
$ws.Cells.Item(33, 11).Value2 = 17000
$ws.Cells.Item(33, 12).Value2 = 18000
$ws.Cells.Item(33, 13).Value2 = 17500
$ws.Cells.Item(33, 16).Value2 = 972
$ws.Cells.Item(34, 4).Value2 = 44630
$ws.Cells.Item(34, 10).Value2 = 60
$ws.Cells.Item(34, 11).Value2 = 15000
$ws.Cells.Item(34, 12).Value2 = 16000
$ws.Cells.Item(34, 13).Value2 = 15500
$ws.Cells.Item(34, 16).Value2 = 861
$ws.Cells.Item(35, 4).Value2 = 44656
$ws.Cells.Item(36, 4).Value2 = 44763
$ws.Cells.Item(36, 10).Value2 = 80
$ws.Cells.Item(36, 11).Value2 = 17000
$ws.Cells.Item(36, 12).Value2 = 18000
$ws.Cells.Item(36, 13).Value2 = 17500
$ws.Cells.Item(36, 16).Value2 = 972
$ws.Cells.Item(37, 4).Value2 = 44754
$ws.Cells.Item(37, 10).Value2 = 80
$ws.Cells.Item(37, 11).Value2 = 16000
$ws.Cells.Item(37, 13).Value2 = 16500
$ws.Cells.Item(37, 16).Value2 = 917
